# Scaling updates to better match EMEP, etc.
# Apply changes to the "year" worksheet of the EMEP_NFR14_scaling_mapping workbook:
#   - Row 2 ("mkd"/"all"): start_scaling_year 1990 -> 2000; add a duplicated
#     "select_scaling_year" header in H1 plus a new "Comment" header in I1;
#     add a "NA" value in H2 and a new comment in I2 explaining the change.
#   - New row 3 ("fin"/"all"): a brand-new scaling-year exception entry
#     (1982-2020) with its own comment.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("year")

# --- Header row (row 1): extend with a duplicated select_scaling_year
#     column and a trailing Comment column ---
$ws.Range("H1").Value = "select_scaling_year"
$ws.Range("I1").Value = "Comment"

# --- Existing data row (row 2): update scaling start year ---
$ws.Range("F2").Value = 2000
$ws.Range("H2").Value = "NA"

# --- New data row (row 3): new exception for Finland ---
$ws.Range("A3").Value = "fin"
$ws.Range("B3").Value = "all"
$ws.Range("C3").Value = "NA"
$ws.Range("D3").Value = "NA"
$ws.Range("E3").Value = "NA"
$ws.Range("F3").Value = 1982
$ws.Range("G3").Value = 2020
$ws.Range("H3").Value = "NA"

# --- Comments (written after the "fin" label so new shared strings keep the
#     same append order Excel itself would produce) ---
$ws.Range("I2").Value = "Scale from 2000 so as to be closer to EMEP trend"
$ws.Range("I3").Value = "Don't scale 1981 to avoid reporting inconsistency in inventory"

# --- Column F width (used by the new selection/layout) ---
$ws.Columns.Item(6).ColumnWidth = 11.75

# --- View state: mirror the saved selections left on the other two sheets
#     while editing, then come back to "year" (the tab that stays active). ---
$wsMap = $wb.Worksheets.Item("map")
$null = $wsMap.Range("B37").Select()

$wsMethod = $wb.Worksheets.Item("method")
$null = $wsMethod.Range("C35").Select()

$null = $ws.Activate()
$null = $ws.Range("A3:XFD3").Select()

$wb.Save()
